$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.762.52"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.633.96"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.77"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.55"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.27"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "1.859.33"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "1.633.37"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.558"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.13"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "25.784.79"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.46"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.49"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.97"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.76"
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.49"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.906"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Value = "1.131.77"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.58"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.81"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.799"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "1.768.32"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.46"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0509"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.416"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("E51").Value = "  +2.86%  "
